$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws.Name = "AccessQuery"
$ws.Columns.Item(1).ColumnWidth = 17.67
$ws.Columns.Item(2).ColumnWidth = 11.5
Write-Output "done"
